$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) stays text-typed like the source data (avoids Excel
# auto-converting numeric-looking strings such as "243.47" into real numbers).
$ws.Range("D2:D51").NumberFormat = "@"

# Update price (D) and volume (E) values for rows 2-40 and 45-51
$ws.Range('D2').Value = '36.249.68'
$ws.Range('E2').Value = '  +1.61%  '
$ws.Range('D3').Value = '2.019.72'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '243.47'
$ws.Range('E5').Value = '  -1.38%  '
$ws.Range('D6').Value = '0.653'
$ws.Range('E6').Value = '  -5.79%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '43.76'
$ws.Range('E8').Value = '  +1.44%  '
$ws.Range('D9').Value = '60.23'
$ws.Range('E9').Value = '  +5.50%  '
$ws.Range('D10').Value = '0.355'
$ws.Range('E10').Value = '  -0.67%  '
$ws.Range('D11').Value = '0.0709'
$ws.Range('E11').Value = '  -6.10%  '
$ws.Range('D12').Value = '0.0980'
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('D13').Value = '14.09'
$ws.Range('E13').Value = '  -3.27%  '
$ws.Range('D14').Value = '2.315.71'
$ws.Range('E14').Value = '  +6.46%  '
$ws.Range('D15').Value = '0.793'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').Value = '2.012.65'
$ws.Range('E16').Value = '  +6.08%  '
$ws.Range('D17').Value = '4.82'
$ws.Range('E17').Value = '  -4.54%  '
$ws.Range('D18').Value = '36.230.99'
$ws.Range('E18').Value = '  +1.56%  '
$ws.Range('D19').Value = '70.49'
$ws.Range('E19').Value = '  -4.36%  '
$ws.Range('D20').Value = '0.0₃0802'
$ws.Range('E20').Value = '  -3.61%  '
$ws.Range('D21').Value = '233.75'
$ws.Range('E21').Value = '  -5.28%  '
$ws.Range('D22').Value = '12.50'
$ws.Range('E22').Value = '  -4.04%  '
$ws.Range('D23').Value = '4.83'
$ws.Range('E23').Value = '  -6.99%  '
$ws.Range('D24').Value = '0.999'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E25').Value = '  -9.36%  '
$ws.Range('D26').Value = '167.74'
$ws.Range('D27').Value = '8.60'
$ws.Range('E27').Value = '  -0.98%  '
$ws.Range('D28').Value = '19.55'
$ws.Range('E28').Value = '  +6.37%  '
$ws.Range('D29').Value = '1.90'
$ws.Range('E29').Value = '  -11.50%  '
$ws.Range('E30').Value = '  -6.42%  '
$ws.Range('D31').Value = '21.30'
$ws.Range('E31').Value = '  +49.89%  '
$ws.Range('D32').Value = '4.28'
$ws.Range('E32').Value = '  -2.35%  '
$ws.Range('D33').Value = '0.0572'
$ws.Range('E33').Value = '  -6.02%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('E35').Value = '  +2.17%  '
$ws.Range('D36').Value = '0.0868'
$ws.Range('E36').Value = '  +18.40%  '
$ws.Range('D37').Value = '3.93'
$ws.Range('E37').Value = '  -7.92%  '
$ws.Range('D38').Value = '2.11'
$ws.Range('E38').Value = '  +7.45%  '
$ws.Range('D39').Value = '0.839'
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('E40').Value = '  -12.12%  '
$ws.Range('D45').Value = '15.25'
$ws.Range('E45').Value = '  -11.07%  '
$ws.Range('D46').Value = '1.300.10'
$ws.Range('E46').Value = '  -1.54%  '
$ws.Range('D47').Value = '0.0813'
$ws.Range('E47').Value = '  +0.02%  '
$ws.Range('D48').Value = '2.77'
$ws.Range('E48').Value = '  +1.11%  '
$ws.Range('D49').Value = '2.238.47'
$ws.Range('E49').Value = '  +7.89%  '
$ws.Range('D50').Value = '2.17'
$ws.Range('E50').Value = '  -8.19%  '
$ws.Range('D51').Value = '3.74'
$ws.Range('E51').Value = '  +12.76%  '

# Rows 41-42 swap: VeChain/Aave -> Aave/VeChain (with updated price/volume)
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '94.79'
$ws.Range('E41').Value = '  -4.53%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '0.0211'
$ws.Range('E42').Value = '  -8.15%  '

# Rows 43-44 swap: ARBITRUM/HuobiToken -> HuobiToken/ARBITRUM (with updated price/volume)
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').Value = '2.79'
$ws.Range('E43').Value = '  +16.19%  '

$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D44').Value = '1.09'
$ws.Range('E44').Value = '  +0.33%  '
